$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "L2" = "right"
    "L3" = "center"
    "L4" = "left"
    "L5" = "center"
    "L6" = "right"
    "L7" = "left"
    "L8" = "center"
    "L9" = "right"
    "L10" = "left"
    "L11" = "center"
    "L12" = "right"
    "L13" = "left"
    "L14" = "center"
    "L15" = "left"
    "L16" = "right"
    "L17" = "left"
    "D18" = "book//book_40.jpg"
    "L18" = "center"
    "L19" = "right"
    "L20" = "center"
    "L21" = "left"
    "L22" = "right"
    "L23" = "right"
    "L24" = "center"
    "L25" = "left"
    "L26" = "center"
    "L27" = "right"
    "L28" = "left"
    "L29" = "left"
    "L30" = "center"
    "L31" = "right"
    "D32" = "book//book_08.jpg"
    "L32" = "center"
    "L33" = "right"
    "L34" = "left"
    "L35" = "left"
    "L36" = "right"
    "L37" = "center"
    "L38" = "center"
    "L39" = "right"
    "L40" = "left"
    "L41" = "right"
    "L42" = "left"
    "L43" = "center"
    "L44" = "right"
    "L45" = "center"
    "L46" = "left"
    "L47" = "left"
    "L48" = "right"
    "L49" = "center"
    "L50" = "right"
    "L51" = "left"
    "L52" = "center"
    "L53" = "right"
    "L54" = "left"
    "L55" = "center"
    "L56" = "left"
    "L57" = "center"
    "L58" = "right"
    "L59" = "right"
    "L60" = "center"
    "L61" = "left"
    "D62" = "book//book_35.jpg"
    "L62" = "center"
    "L63" = "left"
    "L64" = "right"
    "D65" = "book//book_38.jpg"
    "L65" = "center"
    "D66" = "book//book_18.jpg"
    "L66" = "left"
    "L67" = "right"
    "L68" = "left"
    "L69" = "center"
    "L70" = "right"
    "L71" = "right"
    "L72" = "center"
    "L73" = "left"
    "L74" = "left"
    "L75" = "center"
    "L76" = "right"
    "L77" = "center"
    "D78" = "book//book_32.jpg"
    "L78" = "left"
    "L79" = "right"
    "L80" = "right"
    "L81" = "center"
    "L82" = "left"
    "L83" = "right"
    "L84" = "left"
    "L85" = "center"
    "L86" = "left"
    "L87" = "center"
    "L88" = "right"
    "L89" = "left"
    "L90" = "right"
    "L91" = "center"
    "L92" = "left"
    "L93" = "center"
    "L94" = "right"
    "L95" = "left"
    "D96" = "book//book_36.jpg"
    "L96" = "center"
    "L97" = "right"
    "L98" = "right"
    "L99" = "center"
    "L100" = "left"
    "L101" = "center"
    "L102" = "right"
    "L103" = "left"
    "L104" = "center"
    "L105" = "right"
    "L106" = "left"
    "L107" = "left"
    "L108" = "right"
    "L109" = "center"
    "L110" = "right"
    "L111" = "center"
    "L112" = "left"
    "L113" = "right"
    "L114" = "left"
    "L115" = "center"
    "L116" = "center"
    "L117" = "right"
    "L118" = "left"
    "L119" = "left"
    "D120" = "book//book_25.jpg"
    "L120" = "center"
    "L121" = "right"
    "L122" = "left"
    "L123" = "right"
    "D124" = "book//book_04.jpg"
    "L124" = "center"
    "L125" = "right"
    "L126" = "left"
    "L127" = "center"
    "L128" = "center"
    "L129" = "right"
    "L130" = "left"
    "L131" = "right"
    "L132" = "center"
    "L133" = "left"
    "L134" = "center"
    "L135" = "left"
    "L136" = "right"
    "L137" = "left"
    "L138" = "right"
    "L139" = "center"
    "L140" = "left"
    "L141" = "center"
    "L142" = "right"
    "L143" = "left"
    "L144" = "right"
    "L145" = "center"
    "L146" = "center"
    "L147" = "left"
    "L148" = "right"
    "L149" = "center"
    "L150" = "right"
    "L151" = "left"
    "L152" = "left"
    "L153" = "right"
    "L154" = "center"
    "L155" = "right"
    "L156" = "left"
    "L157" = "center"
    "L158" = "right"
    "L159" = "left"
    "L160" = "center"
    "L161" = "center"
    "L162" = "left"
    "L163" = "right"
    "L164" = "center"
    "L165" = "right"
    "L166" = "left"
    "L167" = "left"
    "L168" = "center"
    "L169" = "right"
    "L170" = "left"
    "L171" = "right"
    "L172" = "center"
    "L173" = "center"
    "L174" = "right"
    "D175" = "book//book_25.jpg"
    "L175" = "left"
    "L176" = "right"
    "L177" = "left"
    "L178" = "center"
    "L179" = "center"
    "L180" = "left"
    "L181" = "right"
    "L182" = "right"
    "L183" = "center"
    "L184" = "left"
    "L185" = "right"
    "L186" = "left"
    "L187" = "center"
    "L188" = "right"
    "L189" = "left"
    "L190" = "center"
    "D191" = "book//book_23.jpg"
    "L191" = "right"
    "L192" = "center"
    "L193" = "left"
    "L194" = "left"
    "L195" = "right"
    "L196" = "center"
    "L197" = "left"
    "L198" = "right"
    "L199" = "center"
    "L200" = "right"
    "L201" = "center"
    "L202" = "left"
    "D203" = "book//book_26.jpg"
    "L203" = "left"
    "L204" = "right"
    "L205" = "center"
    "L206" = "center"
    "L207" = "right"
    "L208" = "left"
    "L209" = "left"
    "L210" = "right"
    "L211" = "center"
    "L212" = "right"
    "L213" = "left"
    "L214" = "center"
    "L215" = "right"
    "L216" = "left"
    "L217" = "center"
    "L218" = "left"
    "L219" = "center"
    "L220" = "right"
    "L221" = "left"
    "L222" = "center"
    "L223" = "right"
    "L224" = "right"
    "L225" = "center"
    "L226" = "left"
    "L227" = "center"
    "L228" = "left"
    "L229" = "right"
    "L230" = "right"
    "L231" = "left"
    "L232" = "center"
    "L233" = "left"
    "L234" = "center"
    "L235" = "right"
    "L236" = "right"
    "L237" = "left"
    "L238" = "center"
    "L239" = "left"
    "L240" = "center"
    "L241" = "right"
    "L242" = "right"
    "L243" = "center"
    "L244" = "left"
    "L245" = "left"
    "D246" = "book//book_17.jpg"
    "L246" = "right"
    "L247" = "center"
    "L248" = "left"
    "L249" = "right"
    "D250" = "book//book_40.jpg"
    "L250" = "center"
    "L251" = "center"
    "L252" = "left"
    "L253" = "right"
    "L254" = "left"
    "L255" = "right"
    "L256" = "center"
    "L257" = "right"
    "L258" = "left"
    "L259" = "center"
    "L260" = "left"
    "L261" = "right"
    "L262" = "center"
    "L263" = "center"
    "L264" = "left"
    "L265" = "right"
    "L266" = "center"
    "L267" = "right"
    "L268" = "left"
    "L269" = "left"
    "L270" = "center"
    "L271" = "right"
    "L272" = "left"
    "D273" = "book//book_05.jpg"
    "L273" = "center"
    "L274" = "right"
    "L275" = "left"
    "L276" = "center"
    "L277" = "right"
    "L278" = "right"
    "L279" = "center"
    "L280" = "left"
    "L281" = "center"
    "L282" = "left"
    "L283" = "right"
    "L284" = "center"
    "D285" = "book//book_34.jpg"
    "L285" = "left"
    "L286" = "right"
    "L287" = "center"
    "L288" = "left"
    "L289" = "right"
    "L290" = "center"
    "L291" = "left"
    "D292" = "book//book_24.jpg"
    "L292" = "right"
    "L293" = "left"
    "D294" = "book//book_07.jpg"
    "L294" = "center"
    "L295" = "right"
    "L296" = "center"
    "L297" = "left"
    "L298" = "right"
    "D299" = "book//book_09.jpg"
    "L299" = "right"
    "L300" = "center"
    "L301" = "left"
    "B302" = "book//book_13.jpg"
    "L302" = "center"
    "A303" = "book//book_13.jpg"
    "L303" = "left"
    "L304" = "right"
    "L305" = "center"
    "L306" = "left"
    "L307" = "right"
    "B308" = "book//book_20.jpg"
    "L308" = "center"
    "A309" = "book//book_20.jpg"
    "L309" = "left"
    "L310" = "right"
    "C311" = "book//book_20.jpg"
    "L311" = "right"
    "C312" = "book//book_20.jpg"
    "L312" = "left"
    "L313" = "center"
    "B314" = "book//book_16.jpg"
    "L314" = "right"
    "A315" = "book//book_16.jpg"
    "L315" = "left"
    "L316" = "center"
    "L317" = "left"
    "L318" = "center"
    "L319" = "right"
    "B320" = "book//book_39.jpg"
    "L320" = "left"
    "A321" = "book//book_39.jpg"
    "L321" = "right"
    "L322" = "center"
    "L323" = "right"
    "L324" = "left"
    "C325" = "book//book_39.jpg"
    "L325" = "center"
    "B326" = "book//book_24.jpg"
    "L326" = "left"
    "A327" = "book//book_24.jpg"
    "L327" = "center"
    "L328" = "right"
    "C329" = "book//book_24.jpg"
    "L329" = "center"
    "L330" = "right"
    "L331" = "left"
    "B332" = "book//book_31.jpg"
    "L332" = "center"
    "A333" = "book//book_31.jpg"
    "L333" = "left"
    "L334" = "right"
    "L335" = "left"
    "L336" = "center"
    "L337" = "right"
    "B338" = "book//book_30.jpg"
    "L338" = "left"
    "A339" = "book//book_30.jpg"
    "L339" = "center"
    "L340" = "right"
    "L341" = "center"
    "C342" = "book//book_30.jpg"
    "L342" = "right"
    "C343" = "book//book_30.jpg"
    "L343" = "left"
    "B344" = "book//book_23.jpg"
    "L344" = "center"
    "A345" = "book//book_23.jpg"
    "L345" = "left"
    "L346" = "right"
    "L347" = "left"
    "L348" = "right"
    "L349" = "center"
    "B350" = "book//book_03.jpg"
    "L350" = "left"
    "A351" = "book//book_03.jpg"
    "L351" = "right"
    "L352" = "center"
    "L353" = "center"
    "C354" = "book//book_03.jpg"
    "L354" = "left"
    "L355" = "right"
    "B356" = "book//book_15.jpg"
    "L356" = "center"
    "A357" = "book//book_15.jpg"
    "L357" = "right"
    "L358" = "left"
    "L359" = "center"
    "L360" = "right"
    "L361" = "left"
}

foreach ($key in $changes.Keys) {
    $ws.Range($key).Value = $changes[$key]
}
